$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 475 (old rows 475..555 shift down to 478..558).
$ws.Range("A475:A477").EntireRow.Insert()

$newRange = $ws.Range("A475:D477")

# Match formatting (font, fill, border, alignment) of the surrounding data rows
# (style index 4 in the original workbook) so the inserted rows look identical
# to the rest of the table instead of taking Excel's default "insert" style.
$newRange.Font.Name = "Calibri"
$newRange.Font.Size = 11
$newRange.Font.Color = 0
$newRange.Interior.Pattern = -4142
$newRange.WrapText = $true
$newRange.VerticalAlignment = -4108
$newRange.HorizontalAlignment = 1
$newRange.Borders.Color = 15132391
$newRange.Borders.LineStyle = 1

# Row 475: new series "a) Anteil nachhaltiger Vergabeverfahren"
$ws.Range("A475").Value = "A_SERIES_123a"
$ws.Range("B475").Value = "K_SERIES"
$ws.Range("C475").Value = "a) Anteil nachhaltiger Vergabeverfahren"
$ws.Range("D475").Value = "XXXa) Anteil nachhaltiger Vergabeverfahren"

# Row 476: new series "b) CO2-Emissionen der Kfz der öffentlichen Hand"
$ws.Range("A476").Value = "A_SERIES_123b"
$ws.Range("B476").Value = "K_SERIES"
$ws.Range("C476").Value = "b) CO2-Emissionen der Kfz der öffentlichen Hand"
$ws.Range("D476").Value = "b) CO2 emissions of vehicles in the public sector"

# Row 477: new series "c) Nachhaltige Textilbeschaffung"
$ws.Range("A477").Value = "A_SERIES_123c"
$ws.Range("B477").Value = "K_SERIES"
$ws.Range("C477").Value = "c) Nachhaltige Textilbeschaffung"
$ws.Range("D477").Value = "XXXNachhaltige Textilbeschaffung"
